$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("osmosys_ai_implementers")

$ws.Range("A7").Value = "c6mzlaxl02wz7lk3k"
$ws.Range("B7").Value = "Plan Internacional"
$ws.Range("C7").Value = "No"
$ws.Range("D7").Value = 20
$ws.Range("E7").Value = "PE"
$ws.Range("F7").Value = 1070058
$ws.Range("G7").Value = "Fundación Plan Ecuador"
$ws.Range("H7").Value = "ACTIVO"

$ws.Range("A8").Value = "c7d1hbal02wz7lk3d"
$ws.Range("B8").Value = "Fundación Tarabita"
$ws.Range("C8").Value = "No"
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = "FT"
$ws.Range("F8").Value = 1070030
$ws.Range("G8").Value = "Fundación Tarabita"
$ws.Range("H8").Value = "ACTIVO"

$ws.Range("A13").Value = "ce7axi1l02wz7lj2b"
$ws.Range("B13").Value = "Servicio Jesuita para los Refugiados"
$ws.Range("C13").Value = "No"
$ws.Range("D13").Value = 21
$ws.Range("E13").Value = "SJR"
$ws.Range("F13").Value = 1126007
$ws.Range("G13").Value = "Servicio Jesuita para Refugiados"
$ws.Range("H13").Value = "ACTIVO"

$ws.Range("A14").Value = "cioc6izl02wz7ll48"
$ws.Range("B14").Value = "ASA"
$ws.Range("C14").Value = "No"
$ws.Range("D14").Value = 9
$ws.Range("E14").Value = "ASA"
$ws.Range("F14").Value = 1070041
$ws.Range("G14").Value = "Asociación Solidaridad y Acción"
$ws.Range("H14").Value = "ACTIVO"

$ws.Range("A15").Value = "cirz10l02wz7lj25"
$ws.Range("B15").Value = "Diálogo Diverso"
$ws.Range("C15").Value = "No"
$ws.Range("D15").Value = 17
$ws.Range("E15").Value = "FDD"
$ws.Range("F15").Value = 1070054
$ws.Range("G15").Value = "Fundación Diálogo Diverso"
$ws.Range("H15").Value = "ACTIVO"

$ws.Range("A16").Value = "cjz7296l02wz7lj28"
$ws.Range("B16").Value = "CARE"
$ws.Range("C16").Value = "Si"
$ws.Range("D16").Value = 25
$ws.Range("E16").Value = "CARE"
$ws.Range("F16").Value = 1274040
$ws.Range("G16").Value = "Cooperativa de asistencia y socorro en todas partes"
$ws.Range("H16").Value = "ACTIVO"

$ws.Range("A17").Value = "ckujt1ll02wz7lk2w"
$ws.Range("B17").Value = "Casa de Acogida Matilde"
$ws.Range("C17").Value = "No"
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = "CAI Matilde"
$ws.Range("F17").Value = 1070040
$ws.Range("G17").Value = "Fundación Casa de Refugio Matilde"
$ws.Range("H17").Value = "ACTIVO"

$ws.Range("A18").Value = "cl06o0ul02wz7lj2f"
$ws.Range("B18").Value = "ALDHEA"
$ws.Range("C18").Value = "No"
$ws.Range("D18").Value = 13
$ws.Range("E18").Value = "ALDHEA"
$ws.Range("F18").Value = 1070050
$ws.Range("G18").Value = "Fundación Alternativas Latinoamericanas de Desarrollo Humano y Estudios Antropológicos"
$ws.Range("H18").Value = "ACTIVO"

$ws.Range("A19").Value = "co9elw4l02wz7lj24"
$ws.Range("B19").Value = "Movimiento de Mujeres de El Oro"
$ws.Range("C19").Value = "No"
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = "MMO"
$ws.Range("F19").Value = 1070051
$ws.Range("G19").Value = "Movimiento de Mujeres de El Oro"
$ws.Range("H19").Value = "ACTIVO"

$ws.Range("A20").Value = "cq2w7rrl02wz7lk3i"
$ws.Range("B20").Value = "Fundación CRISFE"
$ws.Range("C20").Value = "No"
$ws.Range("D20").Value = 18
$ws.Range("E20").Value = "CRISFE"
$ws.Range("F20").Value = 1070056
$ws.Range("G20").Value = "Funsación CRISFE"
$ws.Range("H20").Value = "ACTIVO"

$ws.Range("A21").Value = "csoq5cyl02wz7lj2d"
$ws.Range("B21").Value = "Fundación Alas de Colibrí"
$ws.Range("C21").Value = "No"
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = "FAC"
$ws.Range("F21").Value = 1070047
$ws.Range("G21").Value = "Fundación Alas de Colibrí"
$ws.Range("H21").Value = "ACTIVO"

$ws.Range("A22").Value = "cwg57qll02wz7ll46"
$ws.Range("B22").Value = "Federación de Mujeres de Sucumbíos"
$ws.Range("C22").Value = "No"
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = "FMS"
$ws.Range("F22").Value = 1070027
$ws.Range("G22").Value = "Federación de Mujeres de Sumbíos"
$ws.Range("H22").Value = "ACTIVO"

$ws.Range("A23").Value = "cxwq6z7l02wz7lj2h"
$ws.Range("B23").Value = "World Vision (WV)"
$ws.Range("C23").Value = "No"
$ws.Range("D23").Value = 15
$ws.Range("E23").Value = "WV"
$ws.Range("F23").Value = 1070052
$ws.Range("G23").Value = "Fundacion Vision Mundial"
$ws.Range("H23").Value = "ACTIVO"

$ws.Range("A24").Value = "cy5jaell30hqm563"
$ws.Range("B24").Value = "Fundación Quimera"
$ws.Range("C24").Value = "No"
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = "QUIMERA"
$ws.Range("F24").Value = 1070059
$ws.Range("G24").Value = "Fundación Quimera"
$ws.Range("H24").Value = "ACTIVO"

$ws.Range("A25").Value = "cyf42nl02wz7lk3t"
$ws.Range("B25").Value = "Fundación de las Americas (FUDELA)"
$ws.Range("C25").Value = "No"
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = "FUDELA"
$ws.Range("F25").Value = 1070038
$ws.Range("G25").Value = "Fundación de las Americas para el Desarrollo"
$ws.Range("H25").Value = "ACTIVO"

$ws.Range("A26").Value = "czbjbb4l02wz7lj2c"
$ws.Range("B26").Value = "Alto Comisionado de las Naciones Unidas para los Refugiados (ACNUR)"
$ws.Range("C26").Value = "Si"
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = "ACNUR"
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = "Agencia de la ONU para los Refugiados"
$ws.Range("H26").Value = "ACTIVO"

$ws.Range("A27").Value = "czbut89l02wz7lk3c"
$ws.Range("B27").Value = "FEPP"
$ws.Range("C27").Value = "No"
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = "FEPP"
$ws.Range("F27").Value = 1070008
$ws.Range("G27").Value = "Fondo Ecuatoriano Populorum Progressio"
$ws.Range("H27").Value = "ACTIVO"

$ws.Range("A28").Value = "cdw0vygldkcady79"
$ws.Range("B28").Value = "Corporación e Desarrollo de Ambato y Tungurahua CorpoAmbato"
$ws.Range("C28").Value = "No"
$ws.Range("D28").Value = 27
$ws.Range("E28").Value = "CORPOAMBATO"
$ws.Range("F28").Value = 1070060
$ws.Range("G28").Value = "Corporación de Desarrollo de Ambato y Tungurahua"
$ws.Range("H28").Value = "ACTIVO"

$ws.Range("A29").Value = "cnrhc7dle32cstg2"
$ws.Range("B29").Value = "CORPEI"
$ws.Range("C29").Value = "No"

$ws.Activate()
$ws.Range("B10").Select()
